$wb = $excel.ActiveWorkbook

# The workbook has three sheets: Overview, zh-cn, de-de.
# The handback for file "dfcc41c5-cdd4-4d57-983a-71f1e2646132.md" completed successfully,
# so its status flips from "Ready for handoff" to "Handed back: in sync with en-US"
# and the stale error message is replaced by fresh handback timestamps.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to dfcc41c5-cdd4-4d57-983a-71f1e2646132.md
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 corresponds to dfcc41c5-cdd4-4d57-983a-71f1e2646132.md
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-29 14:56:30"
$wsZhCn.Range("P3").Value = ""

# de-de sheet: row 3 corresponds to dfcc41c5-cdd4-4d57-983a-71f1e2646132.md
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-29 14:56:38"
$wsDeDe.Range("P3").Value = ""

# The Error Detail column (P) no longer needs to be wide since the error text is gone;
# Excel re-auto-fits the column after the content shrinks.
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
